# Team weekly report: add the 3/31/2019 and 4/7/2019 entries after the
# 3/24/2019 entry, moving the trailing "_GoBack" bookmark so it stays
# attached to the new last paragraph (4/7/2019), exactly as it was
# attached to the old last paragraph (3/24/2019) before the edit.

$d = $word.ActiveDocument

$enDash = [char]0x2013

# --- locate the "week ending 3/24/2019" paragraph -------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*week ending 3/24/2019*") {
        $targetPara = $p
        break
    }
}

# The "_GoBack" bookmark currently sits right at the end of that
# paragraph's text (just before its paragraph mark). Drop it now - we'll
# re-create it in the same relative spot once the new paragraphs exist.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Position right after the visible text of the 3/24/2019 paragraph, i.e.
# right before its paragraph mark.
$insertPos = $targetPara.Range.End - 1
$anchor = $d.Range($insertPos, $insertPos)

# Create two fresh, empty paragraphs right after the 3/24/2019 one.
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()

$targetIndex = $targetPara.Index
$para2 = $d.Paragraphs.Item($targetIndex + 1)
$para3 = $d.Paragraphs.Item($targetIndex + 2)

# --- paragraph: week ending 3/31/2019 --------------------------------------
$r = $para2.Range
$r.Collapse(1)
$r.InsertAfter("Team weekly report for week ending 3/31/2019 " + $enDash + " The team decided to meet Thursday evening instead of Sunday as there were some scheduling conflicts. We were able to complete ")

$r.Collapse(0)
$r.InsertAfter("3")

$r.Collapse(0)
$r.InsertAfter(" of the 5 tasks for the 3")

$r.Collapse(0)
$r.InsertAfter("rd")
# Mark just the "rd" run (a real, non-empty range) as superscript. Do NOT
# touch .Font on a collapsed range afterwards - that mutates whatever
# run happens to sit at that boundary instead of the (non-existent) run
# being "typed".
$r.Font.Superscript = $true

$r.Collapse(0)
$r.InsertAfter(" iteration. While developing the generate schedule feature, we had to debug some errors as our unit test began to fail. These errors were resolved Thursday night, but took away an hour of programming time to debug.")

# --- paragraph: week ending 4/7/2019 ---------------------------------------
$r3 = $para3.Range
$r3.Collapse(1)
$r3.InsertAfter("Team weekly report for week ending 4/7/2019 " + $enDash + " The team has some scheduling conflicts for various reasons but we were still able to communicate and progress with items on the iteration.")

# --- re-create the "_GoBack" bookmark --------------------------------------
# It needs to end up collapsed right before the paragraph mark of the new
# last paragraph (para3) - the exact spot where it sat on the old last
# paragraph before the edit. Adding a bookmark directly at that boundary
# position is unreliable, so: insert a filler character after the target
# spot (so the position is no longer the very last one in the paragraph),
# add the bookmark there, then delete the filler again. The bookmark stays
# put, and the trailing boundary position is restored.
$bmPos = $para3.Range.End - 1
$filler = $d.Range($bmPos, $bmPos)
$filler.InsertAfter("X")

$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$fillerRange = $d.Range($bmPos, $bmPos + 1)
$fillerRange.Delete()
